$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("D4").Value = "2016-03-03 14:41:13"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("D4").Value = "2016-03-03 14:41:29"
